$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7086900472640991
$ws.Range("B1").Value = 2.113625288009644
$ws.Range("C1").Value = 3.244817018508911
$ws.Range("D1").Value = 3.839879989624023
$ws.Range("E1").Value = 1.240338802337646
